# Revert "FINFLUX-2918 Nabkisan 3months compunding and subsidy scenarios"
# Restores the original (pre-FINFLUX-2918) numbers on the "Repayment schedule"
# and "Summary" sheets, and the original selection state on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Repayment schedule
# ---------------------------------------------------------------------------
$rs = $wb.Worksheets.Item("Repayment schedule")

$rs.Range("F3").Value = 1634.47
$rs.Range("G3").Value = 8418.86
$rs.Range("K3").Value = 1735
$rs.Range("K3").NumberFormat = "#,##0"
$rs.Range("Q3").Value = 1735
$rs.Range("Q3").NumberFormat = "#,##0"

$rs.Range("F4").Value = 1650.81
$rs.Range("G4").Value = 6768.05
$rs.Range("H4").Value = 84.19

$rs.Range("F5").Value = 1667.32
$rs.Range("G5").Value = 5100.73
$rs.Range("H5").Value = 67.68

$rs.Range("F6").Value = 1683.99
$rs.Range("G6").Value = 3416.74
$rs.Range("H6").Value = 51.01

$rs.Range("F7").Value = 1700.83
$rs.Range("G7").Value = 1715.91
$rs.Range("H7").Value = 34.17

$rs.Range("F8").Value = 1715.91
$rs.Range("H8").Value = 19.09

# restore the previously-selected cell
$rs.Range("L11").Select()

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$sm = $wb.Worksheets.Item("Summary")

$sm.Range("F2").Value = 8337.42

$sm.Range("A3").Value = 356.67
$sm.Range("E3").Value = 356.67
$sm.Range("F3").Value = 337.58

# restore the previously-selected cell
$sm.Range("C8").Select()

# ---------------------------------------------------------------------------
# Re-activate "Repayment schedule" (it was the active/selected tab before the
# edit) and restore its selected cell last, so it remains the active sheet.
# ---------------------------------------------------------------------------
$rs.Activate()
$rs.Range("L11").Select()
